$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column B (the structured table column name follows automatically)
$ws.Range("B1").Value = "ADDITIONAL PRACTICE"

# Fill in the new row 5 with an additional practice date/time entry
$ws.Range("A5").Value = 45934
$ws.Range("B5").Value = "12:30pm - 1:30pm"

# Correct the wording / time range on row 3
$ws.Range("B3").Value = "3:00pm - 11:59pm"

# Move the active selection to match the saved view state
$ws.Range("F7").Select() | Out-Null
